$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# D-column values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the original inline-string cells) instead of auto-converting
# number-looking strings (e.g. "582.88") into numeric values.

$ws.Range("D2").Value = "'62.337.80"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "'2.447.03"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'582.88"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").Value = "'144.21"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.93%  "

$ws.Range("D9").Value = "'2.445.05"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  -3.03%  "

$ws.Range("E11").Value = "  +2.66%  "

$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "  -2.65%  "

$ws.Range("D14").Value = "'26.56"
$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("E15").Value = "  -2.86%  "

$ws.Range("D16").Value = "'2.881.49"

$ws.Range("D17").Value = "'62.200.65"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "'2.441.46"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "'10.94"
$ws.Range("E19").Value = "  -3.07%  "

$ws.Range("D20").Value = "'7.17"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").Value = "'330.72"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("D23").Value = "'1.99"
$ws.Range("E23").Value = "  -3.54%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "'65.92"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "'9.45"
$ws.Range("E26").Value = "  +6.37%  "

$ws.Range("D27").Value = "'625.04"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D29").Value = "'0.0₃0959"
$ws.Range("E29").Value = "  -5.72%  "

$ws.Range("E30").Value = "  -0.29%  "

$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").Value = "'4.94"
$ws.Range("E35").Value = "  -4.55%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("E37").Value = "  -5.18%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").Value = "'150.77"
$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("D40").Value = "'5.31"
$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("D41").Value = "'18.36"
$ws.Range("E41").Value = "  -2.05%  "

$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("D43").Value = "'42.48"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("E45").Value = "  -4.16%  "

$ws.Range("D46").Value = "'143.78"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("D47").Value = "'3.65"
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").Value = "'0.0527"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("D49").Value = "'0.601"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("D50").Value = "'19.61"
$ws.Range("E50").Value = "  -7.00%  "

$ws.Range("D51").Value = "'0.0₆0240"
$ws.Range("E51").Value = "  +10.01%  "
